$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "62.784.16"
Set-TextValue "E2" "  -0.94%  "
Set-TextValue "D3" "3.028.02"
Set-TextValue "E3" "  -1.06%  "
Set-TextValue "E4" "  +0.08%  "
Set-TextValue "D5" "585.57"
Set-TextValue "E5" "  -0.51%  "
Set-TextValue "D6" "148.54"
Set-TextValue "E6" "  -4.48%  "
Set-TextValue "E7" "  +0.09%  "
Set-TextValue "D8" "0.526"
Set-TextValue "E8" "  -2.51%  "
Set-TextValue "D9" "3.026.20"
Set-TextValue "E9" "  -1.08%  "
Set-TextValue "D10" "0.151"
Set-TextValue "E10" "  -3.26%  "
Set-TextValue "E11" "  -0.20%  "
Set-TextValue "D12" "0.444"
Set-TextValue "E12" "  -1.57%  "
Set-TextValue "D13" "0.0000231"
Set-TextValue "E13" "  -2.57%  "
Set-TextValue "D14" "35.23"
Set-TextValue "E14" "  -4.95%  "
Set-TextValue "E15" "  +2.36%  "
Set-TextValue "D16" "3.532.86"
Set-TextValue "E16" "  -0.94%  "
Set-TextValue "D17" "7.06"
Set-TextValue "E17" "  -0.84%  "
Set-TextValue "D18" "62.783.32"
Set-TextValue "E18" "  -1.01%  "
Set-TextValue "D19" "3.030.38"
Set-TextValue "E19" "  -0.91%  "
Set-TextValue "D20" "468.05"
Set-TextValue "E20" "  -1.08%  "
Set-TextValue "D21" "14.00"
Set-TextValue "E21" "  -2.40%  "
Set-TextValue "D22" "0.691"
Set-TextValue "E22" "  -1.95%  "
Set-TextValue "D23" "7.40"
Set-TextValue "E23" "  -1.46%  "
Set-TextValue "D24" "2.35"
Set-TextValue "E24" "  -3.38%  "
Set-TextValue "D25" "80.77"
Set-TextValue "E25" "  +0.12%  "
Set-TextValue "D26" "12.43"
Set-TextValue "E26" "  -3.01%  "
Set-TextValue "D27" "10.34"
Set-TextValue "E27" "  -0.46%  "
Set-TextValue "D28" "0.999"
Set-TextValue "E28" "  +0.12%  "
Set-TextValue "E29" "  +0.19%  "
Set-TextValue "D30" "7.21"
Set-TextValue "E30" "  -3.76%  "
Set-TextValue "E31" "  -0.73%  "
Set-TextValue "D32" "2.15"
Set-TextValue "E32" "  -0.21%  "
Set-TextValue "D33" "27.73"
Set-TextValue "E33" "  +2.36%  "
Set-TextValue "D34" "0.107"
Set-TextValue "E34" "  -4.70%  "
Set-TextValue "E35" "  +0.17%  "
Set-TextValue "D36" "0.0₃0803"
Set-TextValue "E36" "  -2.59%  "
Set-TextValue "D37" "5.78"
Set-TextValue "E37" "  -3.33%  "
Set-TextValue "D38" "2.14"
Set-TextValue "E38" "  -2.79%  "
Set-TextValue "D39" "50.34"
Set-TextValue "E39" "  -0.69%  "
Set-TextValue "D40" "9.02"
Set-TextValue "E40" "  -2.54%  "
Set-TextValue "E41" "  -11.00%  "
Set-TextValue "D42" "426.73"
Set-TextValue "E42" "  -3.56%  "
Set-TextValue "E43" "  +2.06%  "
Set-TextValue "D44" "0.280"
Set-TextValue "E44" "  -2.86%  "
Set-TextValue "D45" "2.804.19"
Set-TextValue "E45" "  +0.45%  "
Set-TextValue "D46" "0.0356"
Set-TextValue "E46" "  -0.57%  "
Set-TextValue "D47" "37.62"
Set-TextValue "E47" "  -8.77%  "
Set-TextValue "D48" "129.00"
Set-TextValue "E48" "  -1.06%  "
Set-TextValue "E49" "  +0.00%  "
Set-TextValue "D50" "24.34"
Set-TextValue "E50" "  -2.74%  "
Set-TextValue "E51" "  -0.45%  "
